$d = $word.ActiveDocument

# Word keeps a "_GoBack" bookmark at the location of the most recent edit.
# Remove the one currently sitting near the end of the document; we will
# drop a fresh one in its new location (right after "10 " in change 2 below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Change 1: " Pentru a realizat comunicarea server-client ma recurs..."
#        -> " Pentru a realiza comunicarea server-client am recurs..."
# ---------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("Pentru a realiza", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posRealiza = $f1.End          # right before the "t" of "realizat"

$f2 = $d.Content
$f2.Find.Execute("server-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posServerDash = $f2.End       # right before "client"

$f3 = $d.Content
$f3.Find.Execute("client ma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posClientMaEnd = $f3.End      # right after "ma"

# Drop temporary bookmarks at each future run boundary so the text on
# either side does not get coalesced back into a single run once edited.
$d.Bookmarks.Add("tmpA", $d.Range($posRealiza, $posRealiza))
$d.Bookmarks.Add("tmpB", $d.Range($posServerDash, $posServerDash))
$d.Bookmarks.Add("tmpC", $d.Range($posClientMaEnd, $posClientMaEnd))

# Apply edits right-to-left so earlier (smaller) offsets stay valid -- plain
# integer offsets are not "live" and won't auto-shift as edits are applied.
$d.Range($posClientMaEnd - 2, $posClientMaEnd).Text = "am"       # "ma" -> "am"
$d.Range($posRealiza, $posRealiza + 1).Delete()                  # drop the "t" -> "realiza"

$d.Bookmarks("tmpA").Delete()
$d.Bookmarks("tmpB").Delete()
$d.Bookmarks("tmpC").Delete()

# ---------------------------------------------------------------------
# Change 2: "...dupa 10secunde de la..." -> "...dupa 10 secunde de la..."
# with the new space as its own run and the _GoBack bookmark left
# collapsed right after it, in front of "secunde".
# ---------------------------------------------------------------------
$f4 = $d.Content
$f4.Find.Execute("10secunde", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $f4.Start + 2      # right between "10" and "secunde"

# Insert the new space first, then bracket it with a temp bookmark (on its
# left, between "10" and the space) and the real _GoBack bookmark (on its
# right, between the space and "secunde"). Removing the temp bookmark
# afterwards leaves the space as its own run without re-coalescing it.
$d.Range($splitPos, $splitPos).InsertBefore(" ")
$d.Bookmarks.Add("tmpD", $d.Range($splitPos, $splitPos))
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos + 1, $splitPos + 1))
$d.Bookmarks("tmpD").Delete()

